$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new Result column
$ws.Range("E1").Value = "Result"

# Update existing data rows
$ws.Range("C2").Value = 35
$ws.Range("C3").Value = 20

# New rows 6-8: Mickey test rows (establish "Mickey"/"T" shared strings before "PASS")
$ws.Range("B6").Value = "Mickey"
$ws.Range("C6").Value = 100
$ws.Range("D6").Value = "T"

$ws.Range("B7").Value = "Mickey"
$ws.Range("C7").Value = 100
$ws.Range("D7").Value = "T"

$ws.Range("B8").Value = "Mickey"
$ws.Range("C8").Value = 100
$ws.Range("D8").Value = "T"

# Back-fill PASS result for rows 3 and 5
$ws.Range("E3").Value = "PASS"

# New row 5: duplicate of Vaishnavi with updated Age + PASS result
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Vaishnavi"
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "F"
$ws.Range("E5").Value = "PASS"

$ws.Range("E1").Select()
